$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly block (rows 170 and 171),
# pushing the previously existing rows 170-178 down to 172-180.
$ws.Rows.Item(170).Insert()
$ws.Rows.Item(171).Insert()

# New row 170: Ajo Chino Primera, $/caja 10 kilos
$ws.Cells.Item(170,1).Value() = 9
$ws.Cells.Item(170,2).Value() = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(170,3).Value() = "Metropolitana"
$ws.Cells.Item(170,4).Value() = 44578
$ws.Cells.Item(170,5).Value() = 13
$ws.Cells.Item(170,6).Value() = 100112003
$ws.Cells.Item(170,7).Value() = "Ajo"
$ws.Cells.Item(170,8).Value() = "Chino"
$ws.Cells.Item(170,9).Value() = "Primera"
$ws.Cells.Item(170,10).Value() = 520
$ws.Cells.Item(170,11).Value() = 17500
$ws.Cells.Item(170,12).Value() = 18000
$ws.Cells.Item(170,13).Value() = 17750
$ws.Cells.Item(170,14).Value() = "`$/caja 10 kilos"
$ws.Cells.Item(170,15).Value() = "China"
$ws.Cells.Item(170,16).Value() = 1775
$ws.Cells.Item(170,17).Value() = 10
$ws.Cells.Item(170,18).Value() = "Hortaliza"

# New row 171: Ajo Chino Primera, $/malla 10 kilos
$ws.Cells.Item(171,1).Value() = 9
$ws.Cells.Item(171,2).Value() = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(171,3).Value() = "Metropolitana"
$ws.Cells.Item(171,4).Value() = 44578
$ws.Cells.Item(171,5).Value() = 13
$ws.Cells.Item(171,6).Value() = 100112003
$ws.Cells.Item(171,7).Value() = "Ajo"
$ws.Cells.Item(171,8).Value() = "Chino"
$ws.Cells.Item(171,9).Value() = "Primera"
$ws.Cells.Item(171,10).Value() = 340
$ws.Cells.Item(171,11).Value() = 18000
$ws.Cells.Item(171,12).Value() = 18500
$ws.Cells.Item(171,13).Value() = 18250
$ws.Cells.Item(171,14).Value() = "`$/malla 10 kilos"
$ws.Cells.Item(171,15).Value() = "China"
$ws.Cells.Item(171,16).Value() = 1825
$ws.Cells.Item(171,17).Value() = 10
$ws.Cells.Item(171,18).Value() = "Hortaliza"
